# registerUserName.xlsx — "Add files via upload"
#
# Adds a new worksheet "Text_Try_Editor" at the end of the workbook with a
# couple of sample rows of text, and updates the on-screen selection on the
# previously-active sheet ("invalidCredential") so that the new sheet
# becomes the active/selected tab.

$wb = $excel.ActiveWorkbook

# --- Add the new worksheet at the very end of the tab strip ---------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Text_Try_Editor"

# --- Populate its two rows of sample data ----------------------------------
$ws.Range("A1").Value = "invalid_text"
$ws.Range("B1").Value = "valid_text"
$ws.Range("A2").Value = "hello"
$ws.Range("B2").Value = "print""hello"""

# --- Update the selection left on the previously active sheet -------------
$wsInvalid = $wb.Worksheets.Item("invalidCredential")
$wsInvalid.Range("E20").Select()

# --- Make the new sheet active, with B2 selected ---------------------------
$ws.Activate()
$ws.Range("B2").Select()
